$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 26, shifting existing rows 26-53 down to 27-54.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new record's data.
$ws.Range("A26").Value = 5
$ws.Range("B26").Value = 'Macroferia Regional de Talca'
$ws.Range("C26").Value = 'Maule'
$ws.Range("D26").Value = 44467
$ws.Range("E26").Value = 7
$ws.Range("F26").Value = 100112013
$ws.Range("G26").Value = 'Alcachofa'
$ws.Range("H26").Value = 'Madrigal'
$ws.Range("I26").Value = 'Primera'
$ws.Range("J26").Value = 300
$ws.Range("K26").Value = 10000
$ws.Range("L26").Value = 10000
$ws.Range("M26").Value = 10000
$ws.Range("N26").Value = '$/caja 40 unidades'
$ws.Range("O26").Value = 'Provincia del Elquí'
$ws.Range("P26").Value = 250
$ws.Range("Q26").Value = 40
$ws.Range("R26").Value = 'Hortaliza'

# Match the date-formatted style used by the rest of column D.
$ws.Range("D26").NumberFormat = $ws.Range("D27").NumberFormat
